$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add new "04-29" attendance/score columns ---
# Clone Q1's style (bold/centered header format) into the two new header cells,
# then overwrite their text.
$ws.Range("Q1").Copy($ws.Range("R1"))
$ws.Range("Q1").Copy($ws.Range("S1"))
$ws.Range("R1").Value = "04-29_A"
$ws.Range("S1").Value = "04-29_0"

# --- Data rows (2-119) ---
# Column R duplicates column P (the "attendance" flag column: colored style + value).
$ws.Range("P2:P119").Copy($ws.Range("R2"))

# Column S duplicates column Q as it stood before this edit (the text-typed score).
$ws.Range("Q2:Q119").Copy($ws.Range("S2"))

# Column Q itself is retyped from text to a real number, keeping the same value.
for ($r = 2; $r -le 119; $r++) {
    $cell = $ws.Cells.Item($r, 17)
    $txt = $cell.Value2
    if ($txt -ne "") {
        $cell.Value = [double]$txt
    }
}
